$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.010.61"
$ws.Range("E2").Value = "  +2.56%  "

$ws.Range("D3").Value = "1.650.71"
$ws.Range("E3").Value = "  +3.33%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'214.98"
$ws.Range("E5").Value = "  +1.54%  "

$ws.Range("E6").Value = "  +1.50%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +1.56%  "

$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("D10").Value = "'19.87"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("D12").Value = "1.884.50"
$ws.Range("E12").Value = "  +3.42%  "

$ws.Range("D13").Value = "1.657.59"
$ws.Range("E13").Value = "  +3.65%  "

$ws.Range("D14").Value = "'4.08"
$ws.Range("E14").Value = "  +2.03%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'65.32"
$ws.Range("E15").Value = "  +2.80%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.517"
$ws.Range("E16").Value = "  +2.73%  "

$ws.Range("D17").Value = "'240.05"
$ws.Range("E17").Value = "  +4.01%  "

$ws.Range("D18").Value = "27.000.38"
$ws.Range("E18").Value = "  +2.61%  "

$ws.Range("E19").Value = "  +2.72%  "

$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("E22").Value = "  +4.10%  "

$ws.Range("E23").Value = "  +2.54%  "

$ws.Range("D24").Value = "'9.23"
$ws.Range("E24").Value = "  +3.46%  "

$ws.Range("D25").Value = "'145.94"
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  +1.86%  "

$ws.Range("E28").Value = "  +1.81%  "

$ws.Range("D29").Value = "'15.81"
$ws.Range("E29").Value = "  +2.79%  "

$ws.Range("E30").Value = "  +0.50%  "

$ws.Range("E31").Value = "  +1.93%  "

$ws.Range("E32").Value = "  +2.97%  "

$ws.Range("D33").Value = "1.522.27"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("E34").Value = "  +5.22%  "

$ws.Range("E35").Value = "  +8.32%  "

$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("E37").Value = "  +2.43%  "

$ws.Range("E38").Value = "  +3.05%  "

$ws.Range("D39").Value = "'0.884"
$ws.Range("E39").Value = "  +8.08%  "

$ws.Range("E40").Value = "  +2.86%  "

$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").Value = "'2.25"
$ws.Range("E42").Value = "  +4.00%  "

$ws.Range("D43").Value = "'65.66"
$ws.Range("E43").Value = "  +8.34%  "

$ws.Range("D44").Value = "1.790.88"
$ws.Range("E44").Value = "  +3.23%  "

$ws.Range("D45").Value = "'0.773"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("E46").Value = "  -2.03%  "

$ws.Range("D47").Value = "'89.51"

$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("E49").Value = "  +2.62%  "

$ws.Range("E50").Value = "  +1.29%  "

$ws.Range("E51").Value = "  +1.94%  "
